# Y5_B2728_Anatomy_checklist session log update.
#
# - Rename sheet "Checklist" -> "Session".
# - The log window slid forward by one entry: the oldest row (Student ID
#   231249) rolled off the top and the row that used to be last
#   (Student ID 235020) is no longer retained, so both disappear and the
#   remaining five rows shift up to occupy rows 2-6.
# - Every surviving row's Type flips from "Selection" to "Scan", and a
#   few Log Time values tick forward by a second to match their new
#   "Scan" event.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Session"

# Remove the oldest entry (old row 2 / 231249). Everything below shifts
# up one row.
$ws.Rows(2).Delete()

# Remove what is now the trailing row (old row 8 / 235020) - it falls
# outside the trimmed log.
$ws.Rows(7).Delete()

# Column E ("Type"): Selection -> Scan for all five remaining data rows.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Value = "Scan"
}

# Column D ("Log Time") updates that accompany the new Scan events.
$ws.Cells.Item(2, 4).Value = "21:39:55"
$ws.Cells.Item(5, 4).Value = "21:39:56"
$ws.Cells.Item(6, 4).Value = "21:39:57"
